# Apply cryptos.xlsx data refresh (GitHub Actions scheduled update)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "63.725.89"
$ws.Cells.Item(2, 5).Value = "  +0.18%  "
$ws.Cells.Item(3, 4).Value = "3.121.70"
$ws.Cells.Item(3, 5).Value = "  -0.25%  "
$ws.Cells.Item(4, 5).Value = "  -0.03%  "
$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "602.78"
$ws.Cells.Item(5, 5).Value = "  -0.82%  "
$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = "142.68"
$ws.Cells.Item(6, 5).Value = "  -1.99%  "
$ws.Cells.Item(7, 5).Value = "  +0.05%  "
$ws.Cells.Item(8, 4).Value = "3.115.65"
$ws.Cells.Item(8, 5).Value = "  -0.31%  "
$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.518"
$ws.Cells.Item(9, 5).Value = "  -0.52%  "
$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.148"
$ws.Cells.Item(10, 5).Value = "  -1.32%  "
$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "5.32"
$ws.Cells.Item(11, 5).Value = "  +0.04%  "
$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "0.465"
$ws.Cells.Item(12, 5).Value = "  -1.06%  "
$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "0.0000251"
$ws.Cells.Item(13, 5).Value = "  -0.72%  "
$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "35.01"
$ws.Cells.Item(14, 5).Value = "  -0.89%  "
$ws.Cells.Item(15, 4).Value = "3.635.05"
$ws.Cells.Item(15, 5).Value = "  -0.01%  "
$ws.Cells.Item(16, 5).Value = "  +2.67%  "
$ws.Cells.Item(17, 4).Value = "63.798.33"
$ws.Cells.Item(17, 5).Value = "  +0.32%  "
$ws.Cells.Item(18, 4).Value = "3.129.13"
$ws.Cells.Item(18, 5).Value = "  +0.24%  "
$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "6.78"
$ws.Cells.Item(19, 5).Value = "  -0.90%  "
$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "479.06"
$ws.Cells.Item(20, 5).Value = "  +1.12%  "
$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = "14.55"
$ws.Cells.Item(21, 5).Value = "  +0.24%  "
$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "0.702"
$ws.Cells.Item(22, 5).Value = "  -0.89%  "
$ws.Cells.Item(23, 4).NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = "7.58"
$ws.Cells.Item(23, 5).Value = "  -3.79%  "
$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "88.19"
$ws.Cells.Item(24, 5).Value = "  +6.28%  "
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "13.27"
$ws.Cells.Item(25, 5).Value = "  -2.29%  "
$ws.Cells.Item(26, 5).Value = "  -0.05%  "
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "2.73"
$ws.Cells.Item(27, 5).Value = "  -1.99%  "
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "8.20"
$ws.Cells.Item(28, 5).Value = "  -3.09%  "
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "7.00"
$ws.Cells.Item(29, 5).Value = "  -0.17%  "
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "2.04"
$ws.Cells.Item(30, 5).Value = "  -0.92%  "
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "27.13"
$ws.Cells.Item(31, 5).Value = "  +3.88%  "
$ws.Cells.Item(32, 5).Value = "  -0.03%  "
$ws.Cells.Item(33, 5).Value = "  -8.20%  "
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "2.62"
$ws.Cells.Item(34, 5).Value = "  -2.10%  "
$ws.Cells.Item(35, 5).Value = "  -2.53%  "
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "5.98"
$ws.Cells.Item(36, 5).Value = "  +0.30%  "
$ws.Cells.Item(37, 2).Value = "OKB"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "52.53"
$ws.Cells.Item(37, 5).Value = "  -0.04%  "
$ws.Cells.Item(38, 2).Value = "PEPE"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Cells.Item(38, 4).Value = "0.0₃0748"
$ws.Cells.Item(38, 5).Value = "  -3.52%  "
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "2.93"
$ws.Cells.Item(39, 5).Value = "  -2.08%  "
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "434.47"
$ws.Cells.Item(40, 5).Value = "  -4.87%  "
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.0390"
$ws.Cells.Item(41, 5).Value = "  -0.87%  "
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.119"
$ws.Cells.Item(42, 5).Value = "  +0.51%  "
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "8.23"
$ws.Cells.Item(43, 5).Value = "  -0.70%  "
$ws.Cells.Item(44, 4).Value = "2.855.83"
$ws.Cells.Item(44, 5).Value = "  +0.12%  "
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.258"
$ws.Cells.Item(45, 5).Value = "  -3.46%  "
$ws.Cells.Item(46, 2).Value = "Fetch.AI"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.20"
$ws.Cells.Item(46, 5).Value = "  -4.25%  "
$ws.Cells.Item(47, 2).Value = "ThetaToken"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "2.41"
$ws.Cells.Item(47, 5).Value = "  -0.09%  "
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "25.65"
$ws.Cells.Item(49, 5).Value = "  -2.18%  "
$ws.Cells.Item(50, 5).Value = "  -0.52%  "
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "121.36"
$ws.Cells.Item(51, 5).Value = "  +1.94%  "
